# Update cryptos list values (price & volume columns) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.114.11"
Set-TextValue $ws.Range("E2") "  -1.38%  "
Set-TextValue $ws.Range("D3") "3.525.82"
Set-TextValue $ws.Range("E3") "  +0.09%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "585.67"
Set-TextValue $ws.Range("E5") "  -1.14%  "
Set-TextValue $ws.Range("D6") "133.66"
Set-TextValue $ws.Range("E6") "  -0.21%  "
Set-TextValue $ws.Range("D7") "3.525.50"
Set-TextValue $ws.Range("E7") "  +0.15%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.486"
Set-TextValue $ws.Range("E9") "  -1.06%  "
Set-TextValue $ws.Range("E10") "  +0.01%  "
Set-TextValue $ws.Range("E11") "  -0.50%  "
Set-TextValue $ws.Range("E12") "  -2.06%  "
Set-TextValue $ws.Range("D13") "4.125.35"
Set-TextValue $ws.Range("E13") "  +0.10%  "
Set-TextValue $ws.Range("D14") "27.72"
Set-TextValue $ws.Range("E14") "  +0.16%  "
Set-TextValue $ws.Range("D15") "0.119"
Set-TextValue $ws.Range("E15") "  +1.38%  "
Set-TextValue $ws.Range("E16") "  -1.60%  "
Set-TextValue $ws.Range("D17") "3.526.12"
Set-TextValue $ws.Range("E17") "  +0.20%  "
Set-TextValue $ws.Range("D18") "64.121.29"
Set-TextValue $ws.Range("E18") "  -1.34%  "
Set-TextValue $ws.Range("E19") "  -2.90%  "
Set-TextValue $ws.Range("D20") "14.16"
Set-TextValue $ws.Range("E20") "  -2.13%  "
Set-TextValue $ws.Range("D21") "5.64"
Set-TextValue $ws.Range("E21") "  -1.36%  "
Set-TextValue $ws.Range("D22") "384.80"
Set-TextValue $ws.Range("E22") "  -1.93%  "
Set-TextValue $ws.Range("E23") "  -0.85%  "
Set-TextValue $ws.Range("D24") "3.667.01"
Set-TextValue $ws.Range("E24") "  +0.09%  "
Set-TextValue $ws.Range("D25") "73.86"
Set-TextValue $ws.Range("E25") "  -1.41%  "
Set-TextValue $ws.Range("E26") "  +0.07%  "
Set-TextValue $ws.Range("E27") "  +2.76%  "
Set-TextValue $ws.Range("D28") "1.58"
Set-TextValue $ws.Range("E28") "  -0.18%  "
Set-TextValue $ws.Range("E29") "  -2.45%  "
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  -0.06%  "
Set-TextValue $ws.Range("D31") "8.32"
Set-TextValue $ws.Range("E31") "  -0.60%  "
Set-TextValue $ws.Range("E32") "  -1.58%  "
Set-TextValue $ws.Range("D33") "3.536.07"
Set-TextValue $ws.Range("E33") "  +0.33%  "
Set-TextValue $ws.Range("E34") "  -0.02%  "
Set-TextValue $ws.Range("D35") "23.59"
Set-TextValue $ws.Range("E35") "  -2.23%  "
Set-TextValue $ws.Range("E36") "  +1.02%  "
Set-TextValue $ws.Range("D37") "5.39"
Set-TextValue $ws.Range("E37") "  +1.55%  "
Set-TextValue $ws.Range("E38") "  -1.34%  "
Set-TextValue $ws.Range("D39") "6.93"
Set-TextValue $ws.Range("E39") "  -0.76%  "
Set-TextValue $ws.Range("D40") "159.11"
Set-TextValue $ws.Range("E40") "  -5.51%  "
Set-TextValue $ws.Range("E41") "  -2.89%  "
Set-TextValue $ws.Range("D42") "0.815"
Set-TextValue $ws.Range("E42") "  -0.83%  "
Set-TextValue $ws.Range("D43") "26.30"
Set-TextValue $ws.Range("E43") "  +1.74%  "
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  +0.05%  "
Set-TextValue $ws.Range("D45") "41.96"
Set-TextValue $ws.Range("E45") "  -2.23%  "
Set-TextValue $ws.Range("E46") "  -4.21%  "
Set-TextValue $ws.Range("D47") "4.42"
Set-TextValue $ws.Range("E47") "  -0.61%  "
Set-TextValue $ws.Range("E48") "  -2.48%  "
Set-TextValue $ws.Range("D49") "2.462.88"
Set-TextValue $ws.Range("E49") "  +2.06%  "
Set-TextValue $ws.Range("D50") "6.86"
Set-TextValue $ws.Range("E50") "  -1.08%  "
Set-TextValue $ws.Range("D51") "0.913"
Set-TextValue $ws.Range("E51") "  +0.48%  "
